$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Marker_2" column header in K1
$ws.Range("K1").Value = "Marker_2"

# Fill in marker info for rows that had an off-by-one / missing marker entry.
# Rows 13-16 (strain TDY2274 / CNAG_02700) already carry a G418 marker in J;
# they also carry the NAT marker in the new Marker_2 (K) column.
$ws.Range("K13").Value = "NAT"
$ws.Range("K14").Value = "NAT"
$ws.Range("K15").Value = "NAT"
$ws.Range("K16").Value = "NAT"

# Row 31 (strain TDY2202 / CNAG_07901) was missing its marker_1 entry entirely.
$ws.Range("J31").Value = "NAT"

# Leave the active selection on J31, matching where the edits were made.
$ws.Range("J31").Select()
